# Apply the crypto price/volume refresh described by the commit diff.
# Source data is plain text (coinranking scrape), so every write must land
# as a literal string - numeric-looking D-column values get a leading
# apostrophe (Excel's "force text" quote-prefix) so values like 0.0000268
# or 213.14 are not silently reinterpreted/reformatted as numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '87.748.31'
$ws.Range("E2").Value = '  -1.77%  '

$ws.Range("D3").Value = '3.249.44'
$ws.Range("E3").Value = '  -3.75%  '

$ws.Range("E4").Value = '  -0.01%  '

$ws.Range("D5").Value = '''213.14'
$ws.Range("E5").Value = '  -4.28%  '

$ws.Range("D6").Value = '''628.42'
$ws.Range("E6").Value = '  -2.48%  '

$ws.Range("D7").Value = '''0.387'
$ws.Range("E7").Value = '  +14.27%  '

$ws.Range("D8").Value = '''0.713'
$ws.Range("E8").Value = '  +15.48%  '

$ws.Range("E9").Value = '  +0.01%  '

$ws.Range("D10").Value = '3.243.06'
$ws.Range("E10").Value = '  -4.05%  '

$ws.Range("D11").Value = '''0.577'
$ws.Range("E11").Value = '  -4.89%  '

$ws.Range("E12").Value = '  +11.85%  '

$ws.Range("D13").Value = '''0.0000268'
$ws.Range("E13").Value = '  -3.55%  '

$ws.Range("D14").Value = '''34.33'
$ws.Range("E14").Value = '  -2.76%  '

$ws.Range("D15").Value = '''5.48'
$ws.Range("E15").Value = '  +0.43%  '

$ws.Range("D16").Value = '3.868.22'
$ws.Range("E16").Value = '  -2.92%  '

$ws.Range("D17").Value = '87.722.26'
$ws.Range("E17").Value = '  -1.13%  '

$ws.Range("D18").Value = '3.275.53'
$ws.Range("E18").Value = '  -2.27%  '

$ws.Range("D19").Value = '''3.27'
$ws.Range("E19").Value = '  +1.83%  '

$ws.Range("D20").Value = '''14.05'
$ws.Range("E20").Value = '  -4.87%  '

$ws.Range("D21").Value = '''436.40'
$ws.Range("E21").Value = '  -8.03%  '

$ws.Range("D22").Value = '''8.97'
$ws.Range("E22").Value = '  -2.67%  '

$ws.Range("D23").Value = '''5.33'
$ws.Range("E23").Value = '  -2.81%  '

$ws.Range("D24").Value = '''7.39'
$ws.Range("E24").Value = '  -0.52%  '

$ws.Range("D25").Value = '''5.34'
$ws.Range("E25").Value = '  -1.66%  '

$ws.Range("D26").Value = '''12.45'
$ws.Range("E26").Value = '  -9.61%  '

$ws.Range("B27").Value = 'WrappedeETH'
$ws.Range("C27").Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range("D27").Value = '3.444.02'
$ws.Range("E27").Value = '  -1.82%  '

$ws.Range("B28").Value = 'PEPE'
$ws.Range("C28").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D28").Value = '''0.0000142'
$ws.Range("E28").Value = '  +10.12%  '

$ws.Range("D29").Value = '''77.26'
$ws.Range("E29").Value = '  -2.70%  '

$ws.Range("E30").Value = '  +0.05%  '

$ws.Range("D31").Value = '''0.178'
$ws.Range("E31").Value = '  -16.55%  '

$ws.Range("E32").Value = '  +0.73%  '

$ws.Range("B33").Value = 'InternetComputer(DFINITY)'
$ws.Range("C33").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D33").Value = '''8.88'
$ws.Range("E33").Value = '  -5.31%  '

$ws.Range("B34").Value = 'Bittensor'
$ws.Range("C34").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D34").Value = '''569.47'
$ws.Range("E34").Value = '  -5.14%  '

$ws.Range("D35").Value = '''7.28'
$ws.Range("E35").Value = '  +5.34%  '

$ws.Range("E36").Value = '  -10.93%  '

$ws.Range("D37").Value = '''1.97'
$ws.Range("E37").Value = '  -4.82%  '

$ws.Range("D38").Value = '''0.139'
$ws.Range("E38").Value = '  -8.47%  '

$ws.Range("D39").Value = '''22.96'
$ws.Range("E39").Value = '  -5.11%  '

$ws.Range("D40").Value = '''3.29'
$ws.Range("E40").Value = '  +6.83%  '

$ws.Range("D41").Value = '''21.81'
$ws.Range("E41").Value = '  +0.25%  '

$ws.Range("D42").Value = '''0.997'
$ws.Range("E42").Value = '  -0.41%  '

$ws.Range("D43").Value = '''0.402'
$ws.Range("E43").Value = '  -5.01%  '

$ws.Range("D44").Value = '''2.04'
$ws.Range("E44").Value = '  -2.61%  '

$ws.Range("E45").Value = '  +0.05%  '

$ws.Range("D46").Value = '''151.59'
$ws.Range("E46").Value = '  -3.52%  '

$ws.Range("D47").Value = '''0.138'
$ws.Range("E47").Value = '  +21.36%  '

$ws.Range("D48").Value = '''180.07'
$ws.Range("E48").Value = '  -6.77%  '

$ws.Range("B49").Value = 'OKB'
$ws.Range("C49").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D49").Value = '''45.18'
$ws.Range("E49").Value = '  -5.20%  '

$ws.Range("B50").Value = 'ImmutableX'
$ws.Range("C50").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D50").Value = '''1.35'
$ws.Range("E50").Value = '  -2.53%  '

$ws.Range("D51").Value = '''4.26'
$ws.Range("E51").Value = '  -2.11%  '
